# Apply the latest crypto price/volume snapshot onto the coin table (columns D & E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.260.86'
$ws.Cells.Item(2, 5).Value = '  -1.47%  '
$ws.Cells.Item(3, 4).Value = '2.426.18'
$ws.Cells.Item(3, 5).Value = '  -0.98%  '
$ws.Cells.Item(4, 5).Value = '  -0.11%  '
$ws.Cells.Item(5, 4).Value = '''570.84'
$ws.Cells.Item(5, 5).Value = '  -2.03%  '
$ws.Cells.Item(6, 4).Value = '''140.39'
$ws.Cells.Item(6, 5).Value = '  -1.95%  '
$ws.Cells.Item(7, 5).Value = '  +0.15%  '
$ws.Cells.Item(8, 5).Value = '  -0.85%  '
$ws.Cells.Item(9, 4).Value = '2.413.03'
$ws.Cells.Item(9, 5).Value = '  -1.35%  '
$ws.Cells.Item(10, 5).Value = '  -2.82%  '
$ws.Cells.Item(11, 5).Value = '  -0.35%  '
$ws.Cells.Item(12, 5).Value = '  -2.58%  '
$ws.Cells.Item(13, 5).Value = '  -0.94%  '
$ws.Cells.Item(14, 4).Value = '''26.17'
$ws.Cells.Item(14, 5).Value = '  -1.08%  '
$ws.Cells.Item(15, 4).Value = '''0.0000172'
$ws.Cells.Item(15, 5).Value = '  -2.45%  '
$ws.Cells.Item(16, 4).Value = '2.841.39'
$ws.Cells.Item(17, 4).Value = '61.164.24'
$ws.Cells.Item(17, 5).Value = '  -1.63%  '
$ws.Cells.Item(18, 4).Value = '2.408.78'
$ws.Cells.Item(18, 5).Value = '  -1.65%  '
$ws.Cells.Item(19, 4).Value = '''7.76'
$ws.Cells.Item(19, 5).Value = '  +7.30%  '
$ws.Cells.Item(20, 4).Value = '''10.65'
$ws.Cells.Item(20, 5).Value = '  -0.28%  '
$ws.Cells.Item(21, 4).Value = '''323.68'
$ws.Cells.Item(21, 5).Value = '  -0.94%  '
$ws.Cells.Item(22, 4).Value = '''4.07'
$ws.Cells.Item(22, 5).Value = '  -0.78%  '
$ws.Cells.Item(23, 4).Value = '''6.13'
$ws.Cells.Item(23, 5).Value = '  +2.41%  '
$ws.Cells.Item(24, 5).Value = '  +0.06%  '
$ws.Cells.Item(25, 4).Value = '''1.85'
$ws.Cells.Item(25, 5).Value = '  -2.99%  '
$ws.Cells.Item(26, 4).Value = '''64.60'
$ws.Cells.Item(26, 5).Value = '  -1.64%  '
$ws.Cells.Item(27, 4).Value = '''592.16'
$ws.Cells.Item(27, 5).Value = '  -1.11%  '
$ws.Cells.Item(28, 4).Value = '''8.30'
$ws.Cells.Item(28, 5).Value = '  -9.16%  '
$ws.Cells.Item(29, 4).Value = '2.541.49'
$ws.Cells.Item(29, 5).Value = '  -0.90%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0941'
$ws.Cells.Item(30, 5).Value = '  -3.03%  '
$ws.Cells.Item(31, 4).Value = '''7.94'
$ws.Cells.Item(31, 5).Value = '  -0.74%  '
$ws.Cells.Item(32, 5).Value = '  -4.42%  '
$ws.Cells.Item(33, 5).Value = '  -4.11%  '
$ws.Cells.Item(34, 5).Value = '  -1.22%  '
$ws.Cells.Item(35, 5).Value = '  -0.09%  '
$ws.Cells.Item(36, 4).Value = '''1.43'
$ws.Cells.Item(36, 5).Value = '  -0.59%  '
$ws.Cells.Item(37, 4).Value = '''4.63'
$ws.Cells.Item(37, 5).Value = '  -5.27%  '
$ws.Cells.Item(38, 4).Value = '''151.95'
$ws.Cells.Item(38, 5).Value = '  -0.55%  '
$ws.Cells.Item(39, 4).Value = '''0.369'
$ws.Cells.Item(39, 5).Value = '  -2.05%  '
$ws.Cells.Item(40, 4).Value = '''18.23'
$ws.Cells.Item(40, 5).Value = '  -1.07%  '
$ws.Cells.Item(41, 4).Value = '''5.18'
$ws.Cells.Item(41, 5).Value = '  -2.36%  '
$ws.Cells.Item(43, 4).Value = '''1.68'
$ws.Cells.Item(43, 5).Value = '  -2.00%  '
$ws.Cells.Item(44, 4).Value = '''41.29'
$ws.Cells.Item(44, 5).Value = '  -4.42%  '
$ws.Cells.Item(45, 4).Value = '''2.42'
$ws.Cells.Item(45, 5).Value = '  -4.17%  '
$ws.Cells.Item(46, 4).Value = '0.0₆0302'
$ws.Cells.Item(46, 5).Value = '  +8.03%  '
$ws.Cells.Item(47, 4).Value = '''143.28'
$ws.Cells.Item(47, 5).Value = '  +0.94%  '
$ws.Cells.Item(48, 4).Value = '''3.53'
$ws.Cells.Item(48, 5).Value = '  -2.70%  '
$ws.Cells.Item(49, 4).Value = '''0.589'
$ws.Cells.Item(49, 5).Value = '  -2.15%  '
$ws.Cells.Item(50, 4).Value = '''19.62'
$ws.Cells.Item(50, 5).Value = '  -1.25%  '
$ws.Cells.Item(51, 4).Value = '''0.0503'
$ws.Cells.Item(51, 5).Value = '  -3.08%  '
